# Adding 4 search test cases (TestCase_F15..F17 + supporting rows) to the
# "Test Cases" sheet of the F suite workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Test Cases"
$ws2 = $wb.Worksheets.Item(2)   # "Test Case Steps" (used as a formatting donor)

# ---------------------------------------------------------------------------
# 1) Seed formatting for the three new rows by copying existing cell formats
#    BEFORE writing values, so the engine reuses the same style indices the
#    original workbook already has (s=6 / s=2 / s=3 / s=7 ...).
# ---------------------------------------------------------------------------

# Row 16 mirrors row 15's format exactly (A/B/D/E = plain bordered cell,
# C = bordered cell without wrap).
$ws1.Range("A15:E15").Copy() | Out-Null
$ws1.Range("A16:E16").PasteSpecial(-4122) | Out-Null

# Row 18 uses the same pattern as row 15/16 too.
$ws1.Range("A15:E15").Copy() | Out-Null
$ws1.Range("A18:E18").PasteSpecial(-4122) | Out-Null

# Row 17 is the odd one out: A/D/E like normal (s=6), B has no wrap but no
# fill (s=2, same as the C14/C16 style), C has border + wrap but no fill
# (s=3, borrowed from the "Test Case Steps" sheet which already has it).
$ws1.Range("A15").Copy() | Out-Null
$ws1.Range("A17").PasteSpecial(-4122) | Out-Null

$ws1.Range("C14").Copy() | Out-Null
$ws1.Range("B17").PasteSpecial(-4122) | Out-Null

$ws2.Range("A2").Copy() | Out-Null
$ws1.Range("C17").PasteSpecial(-4122) | Out-Null

$ws1.Range("D15:E15").Copy() | Out-Null
$ws1.Range("D17:E17").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Write the new cell values.
# ---------------------------------------------------------------------------

# Row 16 - TestCase_F15
$ws1.Range("A16").Value = "TestCase_F15"
$ws1.Range("B16").Value = "OPQA-226"
$ws1.Range("C16").Value = "Verify that users should be able to select from a list of suggested topics and check selected topic is presented in users type ahead"
$ws1.Range("D16").Value = "Y"
$ws1.Range("E16").Value = "SKIP"

# Row 17 - TestCase_F16
$ws1.Range("A17").Value = "TestCase_F16"
$ws1.Range("B17").Value = "OPQA-231,OPQA-1100"
$ws1.Range("C17").Value = "Verify that Trending now section include articles and posts and able to navigate from tending now section and `nVerify that Maximum count on the trending list is 10"
$ws1.Range("D17").Value = "Y"
$ws1.Range("E17").Value = "SKIP"

# Row 18 - TestCase_F17
$ws1.Range("A18").Value = "TestCase_F17"
$ws1.Range("B18").Value = "OPQA-1098"
$ws1.Range("C18").Value = "Verify that Featured Post is at the top of event stream after login and that feature post should be top in post tab of trending section"
$ws1.Range("D18").Value = "Y"
$ws1.Range("E18").Value = "PASS"

# ---------------------------------------------------------------------------
# 3) Row height for the wrapped row 17 (matches the 30pt height in the diff).
# ---------------------------------------------------------------------------
$ws1.Rows.Item(17).RowHeight = 30

# ---------------------------------------------------------------------------
# 4) Column layout: split the old merged "A:B" column-width group into its
#    own A and B entries, widening B to fit the longer Jira-id strings.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(2).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 5) View state: active cell / selection moves to D17 (matches the diff's
#    <selection activeCell="D17" sqref="D17"/>).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D17").Select() | Out-Null
